$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 1311.5
$ws.Range("I19").Value = 1196.25
$ws.Range("J19").Value = 1369.125
$ws.Range("K19").Value = 1196.25
$ws.Range("L19").Value = 1369.125
$ws.Range("M19").Value = -1021.25
$ws.Range("N19").Value = -1719.125

$ws.Range("H43").Value = 4007.3333
$ws.Range("I43").Value = 2799.7693
$ws.Range("J43").Value = 7147
$ws.Range("K43").Value = 2799.7693
$ws.Range("L43").Value = 7147
$ws.Range("M43").Value = -2730.7693
$ws.Range("N43").Value = -7285

$ws.Range("H74").Value = 100001850
$ws.Range("I74").Value = 100001850
$ws.Range("J74").Value = 0
$ws.Range("K74").Value = 100001850
$ws.Range("L74").Value = 0
$ws.Range("M74").Value = ""
$ws.Range("N74").Value = -100000914

$ws.Range("H77").Value = 100001850
$ws.Range("I77").Value = 100001850
$ws.Range("J77").Value = 0
$ws.Range("K77").Value = 500009250
$ws.Range("L77").Value = 0
$ws.Range("M77").Value = ""
$ws.Range("N77").Value = -500004570

$ws.Range("H88").Value = 28395404
$ws.Range("I88").Value = 166668130
$ws.Range("J88").Value = 3254910.8
$ws.Range("K88").Value = 166668130
$ws.Range("L88").Value = 3254910.8
$ws.Range("M88").Value = -166667724
$ws.Range("N88").Value = -3255722.8

$ws.Range("H91").Value = 28395404
$ws.Range("I91").Value = 166668130
$ws.Range("J91").Value = 3254910.8
$ws.Range("K91").Value = 166668130
$ws.Range("L91").Value = 3254910.8
$ws.Range("M91").Value = -166666726
$ws.Range("N91").Value = -3257718.8

$ws.Range("H131").Value = 1148
$ws.Range("I131").Value = 1169.1428
$ws.Range("J131").Value = 1000
$ws.Range("K131").Value = 3507.4284
$ws.Range("L131").Value = 3000
$ws.Range("M131").Value = 1532.5716
$ws.Range("N131").Value = -13080

$ws.Range("H137").Value = 2930.5
$ws.Range("I137").Value = 2396
$ws.Range("J137").Value = 3999.5
$ws.Range("K137").Value = 7188
$ws.Range("L137").Value = 11998.5
$ws.Range("M137").Value = -4638
$ws.Range("N137").Value = -17098.5

$ws.Range("H138").Value = 5538.6665
$ws.Range("I138").Value = 4807.8184
$ws.Range("J138").Value = 6687.143
$ws.Range("K138").Value = 14423.4552
$ws.Range("L138").Value = 20061.429
$ws.Range("M138").Value = -9283.4552
$ws.Range("N138").Value = -30341.429

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2344.597
$ws.Range("I32").Value = 1204.5156
$ws.Range("J32").Value = 26666.334
$ws.Range("K32").Value = 1204.5156
$ws.Range("L32").Value = 26666.334
$ws.Range("M32").Value = -917.5155999999999
$ws.Range("N32").Value = -27240.334

$ws.Range("H45").Value = 3567.0386
$ws.Range("I45").Value = 3741.0435
$ws.Range("J45").Value = 2233
$ws.Range("K45").Value = 3741.0435
$ws.Range("L45").Value = 2233
$ws.Range("M45").Value = -3364.0435
$ws.Range("N45").Value = -2987

$ws.Range("H61").Value = 90910880
$ws.Range("I61").Value = 100001780
$ws.Range("J61").Value = 1900
$ws.Range("K61").Value = 100001780
$ws.Range("L61").Value = 1900
$ws.Range("M61").Value = -100001568
$ws.Range("N61").Value = -2324

$ws.Range("H63").Value = 3839.8
$ws.Range("I63").Value = 3839.8
$ws.Range("J63").Value = 0
$ws.Range("K63").Value = 3839.8
$ws.Range("L63").Value = 0
$ws.Range("M63").Value = -3153.8
$ws.Range("N63").Value = ""

$ws.Range("H66").Value = 3839.8
$ws.Range("I66").Value = 3839.8
$ws.Range("J66").Value = 0
$ws.Range("K66").Value = 19199
$ws.Range("L66").Value = 0
$ws.Range("M66").Value = -15767
$ws.Range("N66").Value = ""

$ws.Range("H132").Value = 5266001.5
$ws.Range("I132").Value = 6252720.5
$ws.Range("J132").Value = 3498.3333
$ws.Range("K132").Value = 18758161.5
$ws.Range("L132").Value = 10494.9999
$ws.Range("M132").Value = -18755631.5
$ws.Range("N132").Value = -15554.9999

$ws.Range("H133").Value = 92630.5
$ws.Range("I133").Value = 0
$ws.Range("J133").Value = 92630.5
$ws.Range("K133").Value = 0
$ws.Range("L133").Value = 92630.5
$ws.Range("M133").Value = ""
$ws.Range("N133").Value = -97690.5

$ws.Range("H136").Value = 90910880
$ws.Range("I136").Value = 100001780
$ws.Range("J136").Value = 1900
$ws.Range("K136").Value = 300005340
$ws.Range("L136").Value = 5700
$ws.Range("M136").Value = -300002790
$ws.Range("N136").Value = -10800

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H82").Value = 9000
$ws.Range("I82").Value = 9000
$ws.Range("J82").Value = 0
$ws.Range("K82").Value = 9000
$ws.Range("L82").Value = 0
$ws.Range("M82").Value = -8617
$ws.Range("N82").Value = ""

$ws.Range("H85").Value = 9000
$ws.Range("I85").Value = 9000
$ws.Range("J85").Value = 0
$ws.Range("K85").Value = 9000
$ws.Range("L85").Value = 0
$ws.Range("M85").Value = -7674
$ws.Range("N85").Value = ""

$ws.Range("H107").Value = 47169.305
$ws.Range("I107").Value = 3901.1177
$ws.Range("J107").Value = 169762.5
$ws.Range("K107").Value = 3901.1177
$ws.Range("L107").Value = 169762.5
$ws.Range("M107").Value = -1981.1177
$ws.Range("N107").Value = -173602.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3668.653
$ws.Range("I31").Value = 7221.846
$ws.Range("J31").Value = 2385.5557
$ws.Range("K31").Value = 7221.846
$ws.Range("L31").Value = 2385.5557
$ws.Range("M31").Value = -6926.846
$ws.Range("N31").Value = -2975.5557

$ws.Range("H34").Value = 3668.653
$ws.Range("I34").Value = 7221.846
$ws.Range("J34").Value = 2385.5557
$ws.Range("K34").Value = 7221.846
$ws.Range("L34").Value = 2385.5557
$ws.Range("M34").Value = -7019.846
$ws.Range("N34").Value = -2789.5557

$ws.Range("H86").Value = 12271.75
$ws.Range("I86").Value = 10549.111
$ws.Range("J86").Value = 14486.571
$ws.Range("K86").Value = 10549.111
$ws.Range("L86").Value = 14486.571
$ws.Range("M86").Value = -9426.111000000001
$ws.Range("N86").Value = -16732.571

$ws.Range("H89").Value = 12271.75
$ws.Range("I89").Value = 10549.111
$ws.Range("J89").Value = 14486.571
$ws.Range("K89").Value = 52745.55500000001
$ws.Range("L89").Value = 72432.855
$ws.Range("M89").Value = -47129.55500000001
$ws.Range("N89").Value = -83664.855

$ws.Range("H94").Value = 2559.6
$ws.Range("I94").Value = 3000
$ws.Range("J94").Value = 2266
$ws.Range("K94").Value = 3000
$ws.Range("L94").Value = 2266
$ws.Range("M94").Value = -2549
$ws.Range("N94").Value = -3168

$ws.Range("H132").Value = 21280410
$ws.Range("I132").Value = 24393884
$ws.Range("J132").Value = 4992.1665
$ws.Range("K132").Value = 73181652
$ws.Range("L132").Value = 14976.4995
$ws.Range("M132").Value = -73179122
$ws.Range("N132").Value = -20036.4995

$ws.Range("H134").Value = 5684907.5
$ws.Range("I134").Value = 5955331.5
$ws.Range("J134").Value = 6000
$ws.Range("K134").Value = 17865994.5
$ws.Range("L134").Value = 18000
$ws.Range("M134").Value = -17863459.5
$ws.Range("N134").Value = -23070

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 267.85715
$ws.Range("I2").Value = 194
$ws.Range("J2").Value = 349.1
$ws.Range("K2").Value = 1164
$ws.Range("L2").Value = 2094.6
$ws.Range("M2").Value = -1051
$ws.Range("N2").Value = -2320.6

$ws.Range("H57").Value = 16500
$ws.Range("I57").Value = 2500
$ws.Range("J57").Value = 20000
$ws.Range("K57").Value = 7500
$ws.Range("L57").Value = 60000
$ws.Range("M57").Value = -6941
$ws.Range("N57").Value = -61118

$ws.Range("H68").Value = 7839.54
$ws.Range("I68").Value = 3866.3333
$ws.Range("J68").Value = 8093.149
$ws.Range("K68").Value = 11598.9999
$ws.Range("L68").Value = 24279.447
$ws.Range("M68").Value = -10787.9999
$ws.Range("N68").Value = -25901.447

$ws.Range("H71").Value = 7839.54
$ws.Range("I71").Value = 3866.3333
$ws.Range("J71").Value = 8093.149
$ws.Range("K71").Value = 34796.9997
$ws.Range("L71").Value = 72838.341
$ws.Range("M71").Value = -30740.9997
$ws.Range("N71").Value = -80950.341

$ws.Range("H132").Value = 10687.125
$ws.Range("I132").Value = 22332.666
$ws.Range("J132").Value = 3699.8
$ws.Range("K132").Value = 200993.994
$ws.Range("L132").Value = 33298.2
$ws.Range("M132").Value = -198463.994
$ws.Range("N132").Value = -38358.2

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 1626.7273
$ws.Range("I80").Value = 1819
$ws.Range("J80").Value = 1396
$ws.Range("K80").Value = 1819
$ws.Range("L80").Value = 1396
$ws.Range("M80").Value = -821
$ws.Range("N80").Value = -3392

$ws.Range("H83").Value = 1626.7273
$ws.Range("I83").Value = 1819
$ws.Range("J83").Value = 1396
$ws.Range("K83").Value = 9095
$ws.Range("L83").Value = 6980
$ws.Range("M83").Value = -4103
$ws.Range("N83").Value = -16964

$ws.Range("H97").Value = 1138.4546
$ws.Range("I97").Value = 932.52
$ws.Range("J97").Value = 1782
$ws.Range("K97").Value = 932.52
$ws.Range("L97").Value = 1782
$ws.Range("M97").Value = -436.52
$ws.Range("N97").Value = -2774

$ws.Range("H132").Value = 5684692
$ws.Range("I132").Value = 7355320
$ws.Range("J132").Value = 4557.8
$ws.Range("K132").Value = 22065960
$ws.Range("L132").Value = 13673.4
$ws.Range("M132").Value = -22063430
$ws.Range("N132").Value = -18733.4

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 5500
$ws.Range("I40").Value = 3500
$ws.Range("J40").Value = 7500
$ws.Range("K40").Value = 3500
$ws.Range("L40").Value = 7500
$ws.Range("M40").Value = -3364
$ws.Range("N40").Value = -7772

$ws.Range("H46").Value = 1609.7
$ws.Range("I46").Value = 1785.5714
$ws.Range("J46").Value = 1199.3334
$ws.Range("K46").Value = 1785.5714
$ws.Range("L46").Value = 1199.3334
$ws.Range("M46").Value = -1597.5714
$ws.Range("N46").Value = -1575.3334

$ws.Range("H132").Value = 7146281
$ws.Range("I132").Value = 10420431
$ws.Range("J132").Value = 2680.6365
$ws.Range("K132").Value = 31261293
$ws.Range("L132").Value = 8041.9095
$ws.Range("M132").Value = -31258763
$ws.Range("N132").Value = -13101.9095

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 417.1111
$ws.Range("I107").Value = 428.52
$ws.Range("J107").Value = 274.5
$ws.Range("K107").Value = 1285.56
$ws.Range("L107").Value = 823.5
$ws.Range("M107").Value = 634.4400000000001
$ws.Range("N107").Value = -4663.5

$ws.Range("H122").Value = 1575.4814
$ws.Range("I122").Value = 1075.9
$ws.Range("J122").Value = 3002.8572
$ws.Range("K122").Value = 3227.7
$ws.Range("L122").Value = 9008.571599999999
$ws.Range("M122").Value = -777.7000000000003
$ws.Range("N122").Value = -13908.5716

$ws.Range("H132").Value = 13519875
$ws.Range("I132").Value = 18520432
$ws.Range("J132").Value = 18369.3
$ws.Range("K132").Value = 55561296
$ws.Range("L132").Value = 55107.89999999999
$ws.Range("M132").Value = -55558766
$ws.Range("N132").Value = -60167.89999999999

$ws.Range("H136").Value = 27779876
$ws.Range("I136").Value = 33335072
$ws.Range("J136").Value = 3900
$ws.Range("K136").Value = 100005216
$ws.Range("L136").Value = 11700
$ws.Range("M136").Value = -100002666
$ws.Range("N136").Value = -16800
